$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row (row 3) for the "engage2020" user, mirroring row 2's
# layout/styling but with the engage-specific values.
$ws.Range("A3").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = "engage2020"
$ws.Range("C3").Value = $ws.Range("C2").Text
$ws.Range("D3").Value = "Attrib N Risk Report"
$ws.Range("E3").Value = $ws.Range("E2").Text
$ws.Range("F3").Value = $ws.Range("F2").Text
$ws.Range("G3").Value = "\\sselvarasuw8\screenshot\chartsNew_engage.png"

# Mirror row 2's formatting onto row 3 (hyperlink style + custom font/wrap).
$ws.Range("A3").Style = $ws.Range("A2").Style.Name
$ws.Range("D3").Font.Color = $ws.Range("D2").Font.Color
$ws.Range("D3").WrapText = $ws.Range("D2").WrapText
$ws.Range("G3").Style = $ws.Range("G2").Style.Name
$ws.Range("G3").WrapText = $ws.Range("G2").WrapText

# Hyperlink on G3 pointing at the new screenshot (mirrors G2's hyperlink).
$ws.Hyperlinks.Add($ws.Range("G3"), "file://///sselvarasuw8/screenshot/chartsNew_engage.png") | Out-Null

# Column B needs to widen slightly to fit the new "engage2020" username.
$ws.Columns.Item(2).ColumnWidth = 11.43

# Move the view/selection to the newly added last cell, matching the author's
# on-screen state after adding the row.
$ws.Range("B1").Select() | Out-Null
$ws.Range("G3").Select() | Out-Null

$wb.Save()
